$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2162966666666667
$ws.Range("H2").Value = 0.6488900000000001
$ws.Range("I2").Value = 0.02888548604596741
$ws.Range("J2").Value = 0.0288854860459674
$ws.Range("M2").Value = 0.303146
$ws.Range("N2").Value = 0.909438
$ws.Range("O2").Value = 0.005142855213700541
$ws.Range("P2").Value = 0.005142855213700542
$ws.Range("Q2").Value = 0.06556946931333334
$ws.Range("R2").Value = 0.5901252238200001
$ws.Range("S2").Value = 0.0001485538725117777
$ws.Range("T2").Value = 0.0001485538725117777
$ws.Range("G3").Value = 0.2162966666666667
$ws.Range("H3").Value = 0.6488900000000001
$ws.Range("I3").Value = 0.02888548604596741
$ws.Range("J3").Value = 0.0288854860459674
$ws.Range("O3").Value = 0.2877784259203595
$ws.Range("P3").Value = 0.2877784259203595
$ws.Range("Q3").Value = 3.669066672761111
$ws.Range("R3").Value = 33.02160005485
$ws.Range("S3").Value = 0.008312619706253008
$ws.Range("T3").Value = 0.008312619706253006
$ws.Range("G4").Value = 0.2162966666666667
$ws.Range("H4").Value = 0.6488900000000001
$ws.Range("I4").Value = 0.02888548604596741
$ws.Range("J4").Value = 0.0288854860459674
$ws.Range("M4").Value = 41.67881
$ws.Range("N4").Value = 125.03643
$ws.Range("O4").Value = 0.7070787188659401
$ws.Range("P4").Value = 0.7070787188659401
$ws.Range("Q4").Value = 9.014987673633334
$ws.Range("R4").Value = 81.13488906270001
$ws.Range("S4").Value = 0.02042431246720262
$ws.Range("T4").Value = 0.02042431246720262
$ws.Range("I5").Value = 0.7789723686414617
$ws.Range("J5").Value = 0.7789723686414615
$ws.Range("M5").Value = 0.303146
$ws.Range("N5").Value = 0.909438
$ws.Range("O5").Value = 0.005142855213700541
$ws.Range("P5").Value = 0.005142855213700542
$ws.Range("Q5").Value = 1.768251527438
$ws.Range("R5").Value = 15.914263746942
$ws.Range("S5").Value = 0.004006142107396401
$ws.Range("T5").Value = 0.004006142107396401
$ws.Range("I6").Value = 0.7789723686414617
$ws.Range("J6").Value = 0.7789723686414615
$ws.Range("O6").Value = 0.2877784259203595
$ws.Range("P6").Value = 0.2877784259203595
$ws.Range("S6").Value = 0.2241714420830938
$ws.Range("T6").Value = 0.2241714420830938
$ws.Range("I7").Value = 0.7789723686414617
$ws.Range("J7").Value = 0.7789723686414615
$ws.Range("M7").Value = 41.67881
$ws.Range("N7").Value = 125.03643
$ws.Range("O7").Value = 0.7070787188659401
$ws.Range("P7").Value = 0.7070787188659401
$ws.Range("Q7").Value = 243.11262376643
$ws.Range("R7").Value = 2188.01361389787
$ws.Range("S7").Value = 0.5507947844509715
$ws.Range("T7").Value = 0.5507947844509714
$ws.Range("G8").Value = 1.438774666666667
$ws.Range("H8").Value = 4.316324
$ws.Range("I8").Value = 0.192142145312571
$ws.Range("J8").Value = 0.192142145312571
$ws.Range("M8").Value = 0.303146
$ws.Range("N8").Value = 0.909438
$ws.Range("O8").Value = 0.005142855213700541
$ws.Range("P8").Value = 0.005142855213700542
$ws.Range("Q8").Value = 0.4361587851013333
$ws.Range("R8").Value = 3.925429065912
$ws.Range("S8").Value = 0.0009881592337923626
$ws.Range("T8").Value = 0.0009881592337923628
$ws.Range("G9").Value = 1.438774666666667
$ws.Range("H9").Value = 4.316324
$ws.Range("I9").Value = 0.192142145312571
$ws.Range("J9").Value = 0.192142145312571
$ws.Range("O9").Value = 0.2877784259203595
$ws.Range("P9").Value = 0.2877784259203595
$ws.Range("Q9").Value = 24.40610972158444
$ws.Range("R9").Value = 219.65498749426
$ws.Range("S9").Value = 0.05529436413101265
$ws.Range("T9").Value = 0.05529436413101265
$ws.Range("G10").Value = 1.438774666666667
$ws.Range("H10").Value = 4.316324
$ws.Range("I10").Value = 0.192142145312571
$ws.Range("J10").Value = 0.192142145312571
$ws.Range("M10").Value = 41.67881
$ws.Range("N10").Value = 125.03643
$ws.Range("O10").Value = 0.7070787188659401
$ws.Range("P10").Value = 0.7070787188659401
$ws.Range("Q10").Value = 59.96641596481333
$ws.Range("R10").Value = 539.6977436833199
$ws.Range("S10").Value = 0.135859621947766
$ws.Range("T10").Value = 0.135859621947766
